# Insert a new data row before the existing row 50 ("Hortaliza, Femacal de La
# Calera - Haba" weekly price sheet), pushing the previous rows 50-110 down to
# 51-111 and populating the newly inserted row 50 with the new price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 50:110 down by one row, creating a blank row 50.
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new observation.
$ws.Cells.Item(50, 1).Value  = 3
$ws.Cells.Item(50, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(50, 3).Value  = "Coquimbo"
$ws.Cells.Item(50, 4).Value  = 44540
$ws.Cells.Item(50, 5).Value  = 5
$ws.Cells.Item(50, 6).Value  = 100112026
$ws.Cells.Item(50, 7).Value  = "Haba"
$ws.Cells.Item(50, 8).Value  = "Sin especificar"
$ws.Cells.Item(50, 9).Value  = "Primera"
$ws.Cells.Item(50, 10).Value = 50
$ws.Cells.Item(50, 11).Value = 8000
$ws.Cells.Item(50, 12).Value = 8000
$ws.Cells.Item(50, 13).Value = 8000
$ws.Cells.Item(50, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(50, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(50, 16).Value = 320
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"

# Note: Rows(...).Insert() already carries the donor row's per-cell styling
# (including the date/time number format used in column D) down into the
# newly created row, so no extra style fix-up is required here.
